$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve text formatting for Price column (D) cells so numeric-looking
# strings like "1.002" or "30.238.86" are not auto-converted to numbers.
$dCells = @("D2","D3","D4","D5","D6","D7","D8","D9","D10","D11","D12","D13","D14","D15","D16","D17","D19","D20","D21","D22","D23","D24","D25","D26","D27","D28","D29","D30","D31","D32","D33","D34","D35","D36","D37","D38","D39","D40","D41","D42","D43","D44","D45","D46","D47","D48","D49","D50","D51")
foreach ($addr in $dCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply updated coin data values
$ws.Range("D2").Value = '30.238.86'
$ws.Range("E2").Value = '  -0.52%  '
$ws.Range("D3").Value = '1.882.28'
$ws.Range("E3").Value = '  -1.19%  '
$ws.Range("D4").Value = '1.002'
$ws.Range("E4").Value = '  +0.38%  '
$ws.Range("D5").Value = '237.54'
$ws.Range("E5").Value = '  -0.62%  '
$ws.Range("D6").Value = '1.002'
$ws.Range("E6").Value = '  +0.30%  '
$ws.Range("D7").Value = '0.4665'
$ws.Range("E7").Value = '  -1.36%  '
$ws.Range("D8").Value = '0.2814'
$ws.Range("E8").Value = '  -0.74%  '
$ws.Range("D9").Value = '0.06536'
$ws.Range("E9").Value = '  -1.93%  '
$ws.Range("D10").Value = '19.56'
$ws.Range("E10").Value = '  +4.54%  '
$ws.Range("D11").Value = '0.07757'
$ws.Range("E11").Value = '  +0.60%  '
$ws.Range("D12").Value = '97.15'
$ws.Range("E12").Value = '  -2.64%  '
$ws.Range("D13").Value = '1.893.20'
$ws.Range("E13").Value = '  -1.09%  '
$ws.Range("D14").Value = '5.094'
$ws.Range("E14").Value = '  -2.00%  '
$ws.Range("D15").Value = '0.6624'
$ws.Range("E15").Value = '  -0.78%  '
$ws.Range("D16").Value = '279.29'
$ws.Range("E16").Value = '  +10.08%  '
$ws.Range("D17").Value = '30.253.42'
$ws.Range("E17").Value = '  -0.47%  '
$ws.Range("E18").Value = '  +0.19%  '
$ws.Range("D19").Value = '2.141.71'
$ws.Range("E19").Value = '  -0.58%  '
$ws.Range("D20").Value = '12.54'
$ws.Range("E20").Value = '  -0.53%  '
$ws.Range("D21").Value = '0.000007265'
$ws.Range("E21").Value = '  -2.24%  '
$ws.Range("D22").Value = '5.329'
$ws.Range("E22").Value = '  -0.61%  '
$ws.Range("D23").Value = '1.003'
$ws.Range("E23").Value = '  +0.42%  '
$ws.Range("D24").Value = '6.146'
$ws.Range("E24").Value = '  -2.68%  '
$ws.Range("D25").Value = '166.59'
$ws.Range("E25").Value = '  -0.44%  '
$ws.Range("D26").Value = '9.248'
$ws.Range("E26").Value = '  -1.77%  '
$ws.Range("D27").Value = '18.94'
$ws.Range("E27").Value = '  +0.49%  '
$ws.Range("D28").Value = '1.984'
$ws.Range("E28").Value = '  -2.95%  '
$ws.Range("D29").Value = '1.374'
$ws.Range("E29").Value = '  +0.37%  '
$ws.Range("D30").Value = '0.09775'
$ws.Range("E30").Value = '  -3.09%  '
$ws.Range("D31").Value = '4.438'
$ws.Range("E31").Value = '  -4.32%  '
$ws.Range("D32").Value = '1.485'
$ws.Range("E32").Value = '  -1.60%  '
$ws.Range("D33").Value = '4.150'
$ws.Range("D34").Value = '0.04685'
$ws.Range("E34").Value = '  -0.52%  '
$ws.Range("D35").Value = '0.7031'
$ws.Range("E35").Value = '  -3.35%  '
$ws.Range("D36").Value = '1.088'
$ws.Range("E36").Value = '  -1.81%  '
$ws.Range("D37").Value = '2.718'
$ws.Range("D38").Value = '0.01860'
$ws.Range("E38").Value = '  -2.64%  '
$ws.Range("D39").Value = '6.681'
$ws.Range("E39").Value = '  +7.27%  '
$ws.Range("D40").Value = '2.516'
$ws.Range("E40").Value = '  -2.93%  '
$ws.Range("D41").Value = '71.90'
$ws.Range("E41").Value = '  -2.41%  '
$ws.Range("D42").Value = '0.8679'
$ws.Range("E42").Value = '  +1.08%  '
$ws.Range("D43").Value = '1.956'
$ws.Range("E43").Value = '  +0.12%  '
$ws.Range("B44").Value = 'PaxDollar'
$ws.Range("C44").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D44").Value = '1.001'
$ws.Range("E44").Value = '  +0.33%  '
$ws.Range("B45").Value = 'Quant'
$ws.Range("C45").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D45").Value = '103.79'
$ws.Range("E45").Value = '  -1.56%  '
$ws.Range("D46").Value = '0.4169'
$ws.Range("E46").Value = '  -1.26%  '
$ws.Range("D47").Value = '983.27'
$ws.Range("E47").Value = '  +0.19%  '
$ws.Range("D48").Value = '7.181'
$ws.Range("E48").Value = '  -2.69%  '
$ws.Range("D49").Value = '9.248'
$ws.Range("E49").Value = '  +5.64%  '
$ws.Range("D50").Value = '0.1158'
$ws.Range("E50").Value = '  -3.07%  '
$ws.Range("D51").Value = '33.88'
$ws.Range("E51").Value = '  -1.95%  '

Write-Host "Updated cryptos list"